# edit.ps1 - applies the CFL15 ESP32 Access Point slide7/slide8 changes
# described in the target diff:
#   1. TextBox 13 on slide 7: grow height, add accent1 outline border
#   2. TextBox 13 on slide 7: bump run size to 12pt, append two paragraphs
#      (one blank, one with the GitHub link to the .ino source file)
#   3. Slide 8: add the ZAS Robotics logo picture (copied from slide 7)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Part 1 & 2: slide 7, TextBox 13 ("Code for Lesson CFL15 (...)")
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$tb = $s7.Shapes.Item(5)

$tr = $tb.TextFrame.TextRange

# Bump every existing run to 12pt (adds sz="1200" to each rPr)
$tr.Font.Size = 12

# Append a blank paragraph, then a paragraph with the full GitHub URL.
# (The new runs inherit the current 12pt / bold / en-GB formatting.)
$tr.InsertAfter("`r`rhttps://github.com/info-zas/zas-robotics-communications/blob/main/01_FoundationProjects/CFL15_ESP32_Access_Point/CFL15_AgricutureField_AP_Node/CFL15_AgricutureField_AP_Node.ino")

# Add the accent1-colored outline border (w="19050" = 1.5pt)
$tb.Line.Visible = $true
$tb.Line.Weight = 1.5
$tb.Line.ForeColor.ObjectThemeColor = 5  # msoThemeColorAccent1

# Grow the box to fit the extra lines (cy 646331 -> 1015663 EMU)
$tb.Left = 66.0
$tb.Top = 417.1855118110236
$tb.Width = 426.0652008503937
$tb.Height = 79.97346456692914

# ---------------------------------------------------------------------
# Part 3: slide 8 - add the ZAS Robotics logo picture (same image/
# position/size as the logo already present on every other slide)
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$logo = $s7.Shapes.Item(4)  # "Picture 7" - the logo on slide 7

$logo.Copy()
$pasted = $s8.Shapes.Paste()
$newLogo = $pasted.Item(1)
$newLogo.Name = "Picture 2"

# Position/size to match the canonical logo placement exactly (EMU):
# off x=9004852 y=1, ext cx=3125372 cy=842386
$newLogo.Left = 709.0434875669292
$newLogo.Top = 0.00007874015748031496
$newLogo.Width = 246.09228346456692
$newLogo.Height = 66.3296092992126
